# Avesh Khan.xlsx - "complate!!-> scrapping whole ipl"
#
# 1. Rename the (only) sheet from "Sheet1" to "Avesh Khan".
# 2. Insert a new first column "matchNo" (value "41st" on the data row),
#    pushing the existing teamName..result columns from A:L to B:M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename sheet -----------------------------------------------------
$ws.Name = "Avesh Khan"

# --- insert new column A, shifting everything else to the right -------
$ws.Columns("A").Insert()

# --- new matchNo column -------------------------------------------------
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "41st"

# --- the rest of the header/data row keep their original text, just one
#     column further to the right now (B:M instead of A:L) -------------
$headers = "teamName", "batterName", "states", "runs", "balls", "fours", "sixes", "sr", "opponentTeamName", "venue", "date", "result"
$data    = "Delhi Capitals", "Avesh Khan", "run out (†Karthik)", "5", "3", "1", "0", "166.66", "Kolkata Knight Riders", "Sharjah", "September 28", "KKR won by 3 wickets (with 10 balls remaining)"

# these data values look like plain numbers - prefix with an apostrophe
# (same as typing '5 into a cell) so Excel keeps storing them as text,
# matching the original "numberStoredAsText" data rather than converting
# them into real numbers.
$numericLooking = @{ 3 = $true; 4 = $true; 5 = $true; 6 = $true; 7 = $true }

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item(1, $col).Value = $headers[$i]

    $value = $data[$i]
    if ($numericLooking.ContainsKey($i)) {
        $value = "'" + $value
    }
    $ws.Cells.Item(2, $col).Value = $value
}
